$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 27.60313666666667
$ws.Range("H2").Value = 82.80941000000001
$ws.Range("I2").Value = 0.6521368039512229
$ws.Range("J2").Value = 0.6521368039512228
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 48.91569033333334
$ws.Range("N2").Value = 146.747071
$ws.Range("O2").Value = 0.6566518775718727
$ws.Range("P2").Value = 0.6566518775718726
$ws.Range("Q2").Value = 1350.226485415346
$ws.Range("R2").Value = 12152.03836873811
$ws.Range("S2").Value = 0.4282268567482908
$ws.Range("T2").Value = 0.4282268567482906

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 27.60313666666667
$ws.Range("H3").Value = 82.80941000000001
$ws.Range("I3").Value = 0.6521368039512229
$ws.Range("J3").Value = 0.6521368039512228
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.9094179999999999
$ws.Range("N3").Value = 2.728254
$ws.Range("O3").Value = 0.01220816946726638
$ws.Range("P3").Value = 0.01220816946726638
$ws.Range("Q3").Value = 25.10278934112667
$ws.Range("R3").Value = 225.92510407014
$ws.Range("S3").Value = 0.007961396618478001
$ws.Range("T3").Value = 0.007961396618477997

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 27.60313666666667
$ws.Range("H4").Value = 82.80941000000001
$ws.Range("I4").Value = 0.6521368039512229
$ws.Range("J4").Value = 0.6521368039512228
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1679836666666667
$ws.Range("N4").Value = 0.503951
$ws.Range("O4").Value = 0.002255039014402017
$ws.Range("P4").Value = 0.002255039014402017
$ws.Range("Q4").Value = 4.636876108767779
$ws.Range("R4").Value = 41.73188497891001
$ws.Range("S4").Value = 0.001470593935637447
$ws.Range("T4").Value = 0.001470593935637447

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 27.60313666666667
$ws.Range("H5").Value = 82.80941000000001
$ws.Range("I5").Value = 0.6521368039512229
$ws.Range("J5").Value = 0.6521368039512228
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 24.49948466666666
$ws.Range("N5").Value = 73.498454
$ws.Range("O5").Value = 0.328884913946459
$ws.Range("P5").Value = 0.328884913946459
$ws.Range("Q5").Value = 676.2626235169045
$ws.Range("R5").Value = 6086.36361165214
$ws.Range("S5").Value = 0.2144779566488167
$ws.Range("T5").Value = 0.2144779566488167

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.768748
$ws.Range("H6").Value = 2.306244
$ws.Range("I6").Value = 0.01816202520090028
$ws.Range("J6").Value = 0.01816202520090028
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 48.91569033333334
$ws.Range("N6").Value = 146.747071
$ws.Range("O6").Value = 0.6566518775718727
$ws.Range("P6").Value = 0.6566518775718726
$ws.Range("Q6").Value = 37.60383911236934
$ws.Range("R6").Value = 338.434552011324
$ws.Range("S6").Value = 0.01192612794867884
$ws.Range("T6").Value = 0.01192612794867883

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.768748
$ws.Range("H7").Value = 2.306244
$ws.Range("I7").Value = 0.01816202520090028
$ws.Range("J7").Value = 0.01816202520090028
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.9094179999999999
$ws.Range("N7").Value = 2.728254
$ws.Range("O7").Value = 0.01220816946726638
$ws.Range("P7").Value = 0.01220816946726638
$ws.Range("Q7").Value = 0.6991132686639999
$ws.Range("R7").Value = 6.292019417975999
$ws.Range("S7").Value = 0.0002217250815213533
$ws.Range("T7").Value = 0.0002217250815213533

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.768748
$ws.Range("H8").Value = 2.306244
$ws.Range("I8").Value = 0.01816202520090028
$ws.Range("J8").Value = 0.01816202520090028
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.1679836666666667
$ws.Range("N8").Value = 0.503951
$ws.Range("O8").Value = 0.002255039014402017
$ws.Range("P8").Value = 0.002255039014402017
$ws.Range("Q8").Value = 0.1291371077826667
$ws.Range("R8").Value = 1.162233970044
$ws.Range("S8").Value = 0.00004095607540858276
$ws.Range("T8").Value = 0.00004095607540858275

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.768748
$ws.Range("H9").Value = 2.306244
$ws.Range("I9").Value = 0.01816202520090028
$ws.Range("J9").Value = 0.01816202520090028
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 24.49948466666666
$ws.Range("N9").Value = 73.498454
$ws.Range("O9").Value = 0.328884913946459
$ws.Range("P9").Value = 0.328884913946459
$ws.Range("Q9").Value = 18.83392983853066
$ws.Range("R9").Value = 169.505368546776
$ws.Range("S9").Value = 0.005973216095291508
$ws.Range("T9").Value = 0.005973216095291507

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.1226433333333333
$ws.Range("H10").Value = 0.36793
$ws.Range("I10").Value = 0.002897505178188969
$ws.Range("J10").Value = 0.002897505178188968
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 48.91569033333334
$ws.Range("N10").Value = 146.747071
$ws.Range("O10").Value = 0.6566518775718727
$ws.Range("P10").Value = 0.6566518775718726
$ws.Range("Q10").Value = 5.999183314781111
$ws.Range("R10").Value = 53.99264983303
$ws.Range("S10").Value = 0.00190265221553201
$ws.Range("T10").Value = 0.001902652215532009

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.1226433333333333
$ws.Range("H11").Value = 0.36793
$ws.Range("I11").Value = 0.002897505178188969
$ws.Range("J11").Value = 0.002897505178188968
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.9094179999999999
$ws.Range("N11").Value = 2.728254
$ws.Range("O11").Value = 0.01220816946726638
$ws.Range("P11").Value = 0.01220816946726638
$ws.Range("Q11").Value = 0.1115340549133333
$ws.Range("R11").Value = 1.00380649422
$ws.Range("S11").Value = 0.0000353732342476128
$ws.Range("T11").Value = 0.00003537323424761279

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.1226433333333333
$ws.Range("H12").Value = 0.36793
$ws.Range("I12").Value = 0.002897505178188969
$ws.Range("J12").Value = 0.002897505178188968
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.1679836666666667
$ws.Range("N12").Value = 0.503951
$ws.Range("O12").Value = 0.002255039014402017
$ws.Range("P12").Value = 0.002255039014402017
$ws.Range("Q12").Value = 0.02060207682555555
$ws.Range("R12").Value = 0.18541869143
$ws.Range("S12").Value = 0.000006533987221247992
$ws.Range("T12").Value = 0.000006533987221247991

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.1226433333333333
$ws.Range("H13").Value = 0.36793
$ws.Range("I13").Value = 0.002897505178188969
$ws.Range("J13").Value = 0.002897505178188968
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 24.49948466666666
$ws.Range("N13").Value = 73.498454
$ws.Range("O13").Value = 0.328884913946459
$ws.Range("P13").Value = 0.328884913946459
$ws.Range("Q13").Value = 3.004698464468889
$ws.Range("R13").Value = 27.04228618022
$ws.Range("S13").Value = 0.0009529457411880982
$ws.Range("T13").Value = 0.0009529457411880981

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 13.83269
$ws.Range("H14").Value = 41.49807
$ws.Range("I14").Value = 0.326803665669688
$ws.Range("J14").Value = 0.3268036656696879
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 48.91569033333334
$ws.Range("N14").Value = 146.747071
$ws.Range("O14").Value = 0.6566518775718727
$ws.Range("P14").Value = 0.6566518775718726
$ws.Range("Q14").Value = 676.6355805169967
$ws.Range("R14").Value = 6089.72022465297
$ws.Range("S14").Value = 0.2145962406593712
$ws.Range("T14").Value = 0.2145962406593711

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 13.83269
$ws.Range("H15").Value = 41.49807
$ws.Range("I15").Value = 0.326803665669688
$ws.Range("J15").Value = 0.3268036656696879
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.9094179999999999
$ws.Range("N15").Value = 2.728254
$ws.Range("O15").Value = 0.01220816946726638
$ws.Range("P15").Value = 0.01220816946726638
$ws.Range("Q15").Value = 12.57969727442
$ws.Range("R15").Value = 113.21727546978
$ws.Range("S15").Value = 0.003989674533019415
$ws.Range("T15").Value = 0.003989674533019414

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 13.83269
$ws.Range("H16").Value = 41.49807
$ws.Range("I16").Value = 0.326803665669688
$ws.Range("J16").Value = 0.3268036656696879
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.1679836666666667
$ws.Range("N16").Value = 0.503951
$ws.Range("O16").Value = 0.002255039014402017
$ws.Range("P16").Value = 0.002255039014402017
$ws.Range("Q16").Value = 2.323665986063333
$ws.Range("R16").Value = 20.91299387457
$ws.Range("S16").Value = 0.0007369550161347395
$ws.Range("T16").Value = 0.0007369550161347393

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 13.83269
$ws.Range("H17").Value = 41.49807
$ws.Range("I17").Value = 0.326803665669688
$ws.Range("J17").Value = 0.3268036656696879
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 24.49948466666666
$ws.Range("N17").Value = 73.498454
$ws.Range("O17").Value = 0.328884913946459
$ws.Range("P17").Value = 0.328884913946459
$ws.Range("Q17").Value = 338.8937765537533
$ws.Range("R17").Value = 3050.04398898378
$ws.Range("S17").Value = 0.1074807954611627
$ws.Range("T17").Value = 0.1074807954611627
